$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.404.76"
$ws.Range("E2").Value = "  +2.21%  "
Set-TextValue $ws.Range("D3") "2.096.88"
$ws.Range("E3").Value = "  -0.08%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.68%  "
Set-TextValue $ws.Range("D5") "343.31"
$ws.Range("E5").Value = "  -0.24%  "
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  -0.54%  "
Set-TextValue $ws.Range("D7") "0.5259"
$ws.Range("E7").Value = "  +1.75%  "
Set-TextValue $ws.Range("D8") "0.4426"
$ws.Range("E8").Value = "  +0.84%  "
Set-TextValue $ws.Range("D9") "54.61"
$ws.Range("E9").Value = "  +3.39%  "
Set-TextValue $ws.Range("D10") "0.09356"
$ws.Range("E10").Value = "  +0.76%  "
Set-TextValue $ws.Range("D11") "1.170"
$ws.Range("E11").Value = "  +0.49%  "
Set-TextValue $ws.Range("D12") "24.78"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "8.575"
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "2.128.47"
$ws.Range("E14").Value = "  +1.25%  "
Set-TextValue $ws.Range("D15") "6.924"
$ws.Range("E15").Value = "  +2.32%  "
Set-TextValue $ws.Range("D16") "101.50"
$ws.Range("E16").Value = "  +1.92%  "
Set-TextValue $ws.Range("D17") "0.00001160"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  -0.59%  "
Set-TextValue $ws.Range("D19") "21.22"
$ws.Range("E19").Value = "  +1.80%  "
Set-TextValue $ws.Range("D20") "0.06689"
$ws.Range("E20").Value = "  +0.56%  "
Set-TextValue $ws.Range("D21") "6.339"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  -0.55%  "
Set-TextValue $ws.Range("D23") "30.437.97"
$ws.Range("E23").Value = "  +2.23%  "
Set-TextValue $ws.Range("D24") "12.57"
$ws.Range("E24").Value = "  +0.61%  "
Set-TextValue $ws.Range("D25") "2.310"
$ws.Range("E25").Value = "  -0.22%  "
Set-TextValue $ws.Range("D26") "21.87"
$ws.Range("E26").Value = "  -0.37%  "
Set-TextValue $ws.Range("D27") "163.25"
$ws.Range("E27").Value = "  +1.13%  "
Set-TextValue $ws.Range("D28") "6.801"
$ws.Range("E28").Value = "  +8.11%  "
Set-TextValue $ws.Range("D29") "2.513"
$ws.Range("E29").Value = "  -0.25%  "
Set-TextValue $ws.Range("D30") "133.58"
$ws.Range("E30").Value = "  +0.34%  "
Set-TextValue $ws.Range("D31") "1.138"
$ws.Range("E31").Value = "  -0.46%  "
Set-TextValue $ws.Range("D32") "0.1050"
$ws.Range("E32").Value = "  +0.09%  "
Set-TextValue $ws.Range("D33") "1.652"
$ws.Range("E33").Value = "  -0.12%  "
Set-TextValue $ws.Range("D34") "6.279"
$ws.Range("E34").Value = "  +1.63%  "
Set-TextValue $ws.Range("D35") "3.870"
$ws.Range("E35").Value = "  -1.70%  "
Set-TextValue $ws.Range("D36") "10.17"
$ws.Range("E36").Value = "  -0.33%  "
Set-TextValue $ws.Range("D37") "0.02639"
$ws.Range("E37").Value = "  +2.19%  "
Set-TextValue $ws.Range("D38") "0.06805"
$ws.Range("E38").Value = "  +1.09%  "
Set-TextValue $ws.Range("D39") "0.7018"
$ws.Range("E39").Value = "  +1.86%  "
Set-TextValue $ws.Range("D40") "12.60"
$ws.Range("E40").Value = "  +0.96%  "
Set-TextValue $ws.Range("D41") "1.344"
$ws.Range("E41").Value = "  +1.86%  "
Set-TextValue $ws.Range("D42") "0.2221"
Set-TextValue $ws.Range("D43") "0.6864"
$ws.Range("E43").Value = "  +1.29%  "
Set-TextValue $ws.Range("D44") "14.45"
$ws.Range("E44").Value = "  +1.27%  "
Set-TextValue $ws.Range("D45") "2.345"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  -0.49%  "
Set-TextValue $ws.Range("D47") "1.386"
$ws.Range("E47").Value = "  +19.48%  "
Set-TextValue $ws.Range("D48") "3.637"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D49") "1.235"
$ws.Range("E49").Value = "  +9.90%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D50") "0.00000000344"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("E51").Value = "  -0.26%  "
